# Update the Area RRP value in the "Template" sheet: cell B2 goes from "66" to "76".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Range("B2").Value = "76"
